# "added multiple invalid login"
# Add a new worksheet "MultipleInvalidLogin" (after the existing "InvalidLogin"
# sheet), populate it with a username/password table (reusing the same
# username/password values that appear in the other sheets, plus a new
# "admin" entry), box every cell with a thin border, auto-fit the password
# column, and make the new sheet the active/selected one.

$wb = $excel.ActiveWorkbook

# Insert the new sheet right after the current last sheet ("InvalidLogin")
# so it lands at the end of the tab strip.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "MultipleInvalidLogin"

# Header row + three sample login attempts.
$newSheet.Range("A1").Value = "Username"
$newSheet.Range("B1").Value = "password"
$newSheet.Range("A2").Value = "abcd"
$newSheet.Range("B2").Value = "xyz"
$newSheet.Range("A3").Value = "admin"
$newSheet.Range("B3").Value = "xyz"
$newSheet.Range("A4").Value = "bhanu"
$newSheet.Range("B4").Value = "pointofsale"

# Box every cell in the table with a thin border on all sides.
$dataRange = $newSheet.Range("A1:B4")
$dataRange.Borders.LineStyle = 1

# Best-fit the password column so "pointofsale" isn't clipped.
$newSheet.Columns.Item(2).AutoFit()

# Make the new sheet the active tab with the same zoom level used by the
# other test-case sheets, and leave the whole table selected.
[void]$newSheet.Activate()
$excel.ActiveWindow.Zoom = 190
[void]$dataRange.Select()
